$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Move the "_GoBack" bookmark from the end of the paragraph that ends with
#    "...imitating the genetic behaviour." to the end of the "Grass growth
#    rate" paragraph.
#
#    Directly adding a collapsed bookmark exactly at the last character
#    position of a paragraph (i.e. right before the paragraph mark) is
#    unreliable in this host, so we insert a temporary placeholder
#    character, wrap the bookmark around it, and then delete the
#    placeholder - leaving a correctly collapsed bookmark behind.
# ---------------------------------------------------------------------------

# Remove the bookmark from its current location (if present).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Find the "Grass growth rate" paragraph and locate the position right
# before its paragraph mark.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Grass growth rate*") {
        $target = $p
    }
}

$insertAt = $d.Range($target.Range.End - 1, $target.Range.End - 1)
$insertAt.InsertAfter("X")
$placeholder = $d.Range($target.Range.End - 2, $target.Range.End - 1)
$placeholder.Bookmarks.Add("_GoBack")
$placeholder.Text = ""

# ---------------------------------------------------------------------------
# 2. After the paragraph ending in "...imitating the genetic behaviour.",
#    insert two new blank paragraphs with the same paragraph formatting as
#    the surrounding body paragraphs (first-line indent + justified + en-GB).
# ---------------------------------------------------------------------------

$genParagraph = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*imitating the genetic behaviour*") {
        $genParagraph = $p
    }
}

$blankParaXml = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:pPr><w:ind w:firstLine='360'/><w:jc w:val='both'/><w:rPr><w:lang w:val='en-GB'/></w:rPr></w:pPr></w:p>"

$insertPoint1 = $d.Range($genParagraph.Range.End, $genParagraph.Range.End)
$null = $insertPoint1.InsertXML($blankParaXml)

$genParagraph2 = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*imitating the genetic behaviour*") {
        $genParagraph2 = $p
    }
}
$insertPoint2 = $d.Range($genParagraph2.Range.End, $genParagraph2.Range.End)
$null = $insertPoint2.InsertXML($blankParaXml)
